$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.615.88'
$ws.Range("E2").Value = '  +4.10%  '
$ws.Range("D3").Value = '3.493.81'
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''582.83'
$ws.Range("E5").Value = '  +3.41%  '
$ws.Range("D6").Value = '''161.36'
$ws.Range("E6").Value = '  +4.51%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '''0.608'
$ws.Range("E8").Value = '  +13.00%  '
$ws.Range("D9").Value = '3.500.39'
$ws.Range("E9").Value = '  +3.12%  '
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("E11").Value = '  +4.14%  '
$ws.Range("D12").Value = '''0.448'
$ws.Range("E12").Value = '  +4.08%  '
$ws.Range("D13").Value = '4.097.03'
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("D14").Value = '''0.135'
$ws.Range("E14").Value = '  +0.75%  '
$ws.Range("E15").Value = '  +4.26%  '
$ws.Range("D16").Value = '''28.73'
$ws.Range("E16").Value = '  +7.39%  '
$ws.Range("D17").Value = '65.628.64'
$ws.Range("E17").Value = '  +3.56%  '
$ws.Range("D18").Value = '3.481.38'
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("E19").Value = '  +4.17%  '
$ws.Range("D20").Value = '''14.35'
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("D21").Value = '''389.61'
$ws.Range("E21").Value = '  +2.04%  '
$ws.Range("D22").Value = '''8.30'
$ws.Range("E22").Value = '  +2.78%  '
$ws.Range("D23").Value = '''0.557'
$ws.Range("E23").Value = '  +4.99%  '
$ws.Range("D24").Value = '''73.45'
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("E26").Value = '  +7.95%  '
$ws.Range("D27").Value = '''10.22'
$ws.Range("E27").Value = '  +9.60%  '
$ws.Range("E28").Value = '  +2.09%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").Value = '''6.33'
$ws.Range("E30").Value = '  +6.61%  '
$ws.Range("D31").Value = '''1.48'
$ws.Range("E31").Value = '  +11.13%  '
$ws.Range("E32").Value = '  +4.32%  '
$ws.Range("D33").Value = '''23.75'
$ws.Range("E33").Value = '  +3.12%  '
$ws.Range("E34").Value = '  +7.65%  '
$ws.Range("E35").Value = '  +10.33%  '
$ws.Range("D36").Value = '''163.20'
$ws.Range("E36").Value = '  +3.16%  '
$ws.Range("E37").Value = '  +7.95%  '
$ws.Range("D38").Value = '3.042.08'
$ws.Range("E38").Value = '  +6.03%  '
$ws.Range("E39").Value = '  +2.71%  '
$ws.Range("D40").Value = '''27.44'
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("D41").Value = '''0.0327'
$ws.Range("E41").Value = '  +4.07%  '
$ws.Range("D42").Value = '''4.60'
$ws.Range("E42").Value = '  +6.51%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''6.55'
$ws.Range("E43").Value = '  +2.99%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''43.09'
$ws.Range("E44").Value = '  +7.31%  '
$ws.Range("D45").Value = '''0.780'
$ws.Range("E45").Value = '  +3.02%  '
$ws.Range("D46").Value = '''25.92'
$ws.Range("E46").Value = '  +11.42%  '
$ws.Range("E47").Value = '  +5.52%  '
$ws.Range("D48").Value = '''320.78'
$ws.Range("E48").Value = '  +12.13%  '
$ws.Range("D49").Value = '''6.77'
$ws.Range("E49").Value = '  +6.38%  '
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("E51").Value = '  +7.37%  '
